$wb = $excel.ActiveWorkbook

# --- Sheet "TestCase_A6": set Runmode column (C) to "Y" for rows 2-3 ---
$wsA6 = $wb.Worksheets.Item("TestCase_A6")
$wsA6.Range("C2").Value = "Y"
$wsA6.Range("C3").Value = "Y"

# --- Sheet "TestCase_A5": set Runmode column (C) to "Y" for rows 2-3 ---
$wsA5 = $wb.Worksheets.Item("TestCase_A5")
$wsA5.Range("C2").Value = "Y"
$wsA5.Range("C3").Value = "Y"
$wsA5.Range("C5").Select()

# --- Sheet "Test Cases": set Runmode column (C) to "Y" for rows 3-7 ---
$wsTestCases = $wb.Worksheets.Item("Test Cases")
$wsTestCases.Range("C3").Value = "Y"
$wsTestCases.Range("C4").Value = "Y"
$wsTestCases.Range("C5").Value = "Y"
$wsTestCases.Range("C6").Value = "Y"
$wsTestCases.Range("C7").Value = "Y"
$wsTestCases.Range("C8").Select()
